# Apply updated crypto market data to Sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.268.43'
$ws.Range('E2').Value = '  +1.40%  '

# Row 3
$ws.Range('D3').Value = '1.806.15'
$ws.Range('E3').Value = '  +3.01%  '

# Row 4
$ws.Range('E4').Value = '  -0.26%  '

# Row 5
$ws.Range('D5').Value = '''336.33'
$ws.Range('E5').Value = '  +0.26%  '

# Row 6
$ws.Range('E6').Value = '  -0.20%  '

# Row 7
$ws.Range('D7').Value = '''0.4617'
$ws.Range('E7').Value = '  +20.67%  '

# Row 8
$ws.Range('D8').Value = '''0.3718'
$ws.Range('E8').Value = '  +9.20%  '

# Row 9
$ws.Range('D9').Value = '''45.12'
$ws.Range('E9').Value = '  -2.64%  '

# Row 10
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '''0.07662'
$ws.Range('E10').Value = '  +6.07%  '

# Row 11
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '''1.153'
$ws.Range('E11').Value = '  +3.49%  '

# Row 12
$ws.Range('D12').Value = '''22.40'
$ws.Range('E12').Value = '  -0.48%  '

# Row 13
$ws.Range('E13').Value = '  -0.25%  '

# Row 14
$ws.Range('D14').Value = '''6.345'
$ws.Range('E14').Value = '  +2.89%  '

# Row 15
$ws.Range('D15').Value = '''7.489'
$ws.Range('E15').Value = '  +4.74%  '

# Row 16
$ws.Range('D16').Value = '1.806.09'
$ws.Range('E16').Value = '  +2.85%  '

# Row 17
$ws.Range('E17').Value = '  +3.87%  '

# Row 18
$ws.Range('D18').Value = '''0.06721'
$ws.Range('E18').Value = '  +1.76%  '

# Row 19
$ws.Range('D19').Value = '''81.94'
$ws.Range('E19').Value = '  +3.72%  '

# Row 20
$ws.Range('E20').Value = '  -0.22%  '

# Row 21
$ws.Range('D21').Value = '''17.49'
$ws.Range('E21').Value = '  +4.62%  '

# Row 22
$ws.Range('D22').Value = '''6.419'
$ws.Range('E22').Value = '  +3.08%  '

# Row 23
$ws.Range('D23').Value = '28.251.53'
$ws.Range('E23').Value = '  +1.26%  '

# Row 24
$ws.Range('D24').Value = '''11.88'
$ws.Range('E24').Value = '  +1.86%  '

# Row 25
$ws.Range('D25').Value = '''2.405'
$ws.Range('E25').Value = '  +0.80%  '

# Row 26
$ws.Range('D26').Value = '''20.80'
$ws.Range('E26').Value = '  +4.61%  '

# Row 27
$ws.Range('D27').Value = '''153.58'
$ws.Range('E27').Value = '  +0.70%  '

# Row 28
$ws.Range('D28').Value = '''2.379'
$ws.Range('E28').Value = '  +2.71%  '

# Row 29
$ws.Range('D29').Value = '2.011.43'
$ws.Range('E29').Value = '  +2.77%  '

# Row 30
$ws.Range('D30').Value = '''133.34'
$ws.Range('E30').Value = '  +0.89%  '

# Row 31
$ws.Range('D31').Value = '''1.258'
$ws.Range('E31').Value = '  -1.14%  '

# Row 32
$ws.Range('D32').Value = '''4.031'
$ws.Range('E32').Value = '  +0.17%  '

# Row 33
$ws.Range('D33').Value = '''0.09571'
$ws.Range('E33').Value = '  +8.65%  '

# Row 34
$ws.Range('D34').Value = '''5.867'
$ws.Range('E34').Value = '  +0.48%  '

# Row 35
$ws.Range('D35').Value = '''0.2222'
$ws.Range('E35').Value = '  +5.57%  '

# Row 36
$ws.Range('D36').Value = '''12.13'
$ws.Range('E36').Value = '  -0.71%  '

# Row 37
$ws.Range('D37').Value = '''0.06368'
$ws.Range('E37').Value = '  +3.38%  '

# Row 38
$ws.Range('D38').Value = '''0.02354'
$ws.Range('E38').Value = '  +3.01%  '

# Row 39
$ws.Range('D39').Value = '''5.263'
$ws.Range('E39').Value = '  +2.39%  '

# Row 40
$ws.Range('D40').Value = '''0.6650'
$ws.Range('E40').Value = '  +0.94%  '

# Row 41
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '''1.509'
$ws.Range('E41').Value = '  +0.42%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.238'
$ws.Range('E42').Value = '  +2.48%  '

# Row 43
$ws.Range('D43').Value = '''8.249'
$ws.Range('E43').Value = '  +3.22%  '

# Row 44
$ws.Range('D44').Value = '''14.30'
$ws.Range('E44').Value = '  +3.66%  '

# Row 45
$ws.Range('D45').Value = '''0.9982'
$ws.Range('E45').Value = '  -0.27%  '

# Row 46
$ws.Range('D46').Value = '''0.6134'
$ws.Range('E46').Value = '  +0.95%  '

# Row 47
$ws.Range('D47').Value = '''3.824'
$ws.Range('E47').Value = '  +0.01%  '

# Row 48
$ws.Range('D48').Value = '''129.98'
$ws.Range('E48').Value = '  +2.93%  '

# Row 49
$ws.Range('E49').Value = '  +2.37%  '

# Row 50
$ws.Range('D50').Value = '''0.07160'
$ws.Range('E50').Value = '  +2.63%  '

# Row 51
$ws.Range('D51').Value = '''1.178'
$ws.Range('E51').Value = '  +0.52%  '
